$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.558.74'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.882.52'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '249.42'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4745'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2934'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.78%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06530'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.06'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07732'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7470'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '96.85'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.881.43'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.280'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '276.13'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.546.40'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.25'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -3.08%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('B19').Style = "Normal"
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C19').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007544'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'Dai'
$ws.Range('B20').Style = "Normal"
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C20').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.128.07'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.322'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.246'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('B25').Style = "Normal"
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C25').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.216'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Monero'
$ws.Range('B26').Style = "Normal"
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C26').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.40'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.88'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.915'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('B29').Style = "Normal"
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C29').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.344'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Stellar'
$ws.Range('B30').Style = "Normal"
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C30').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09711'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.81%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.504'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.305'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.158'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04864'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.06%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.128'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6985'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01897'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.773'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.311'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '74.36'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +5.85%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.022'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4246'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.50%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('B44').Style = "Normal"
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C44').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8404'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '102.28'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.479'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Elrond'
$ws.Range('B48').Style = "Normal"
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('C48').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '35.72'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Aptos'
$ws.Range('B49').Style = "Normal"
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C49').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.045'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '916.40'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05724'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.73%  '
$ws.Range('E51').Style = "Normal"
